$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that originally sat right after the
#    inline picture (it will be re-created later, anchored after "140").
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) The paragraph "(x, y) = (623.5918, 322.0938)" gets an explicit
#    "no numbering" numPr (ilvl=0/numId=0) plus a first-line indent of
#    100 characters (240 twips @ 1 char = 2.4*10pt households -> Word
#    stores both the twip value and the character-unit value).
# ---------------------------------------------------------------------------
$targetRng = $d.Content.Duplicate
$targetRng.Find.Execute("(x, y) = (623.5918, 322.0938)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$targetPara = $targetRng.Paragraphs(1)
$targetPara.Range.ListFormat.RemoveNumbers()
$targetPara.Range.ParagraphFormat.CharacterUnitFirstLineIndent = 100
$targetPara.Range.ParagraphFormat.FirstLineIndent = 12

# ---------------------------------------------------------------------------
# 3) Split the bold run "iteration 161" into "iteration " + "140", and move
#    the "_GoBack" bookmark so that it now sits right after the new "140"
#    run.
# ---------------------------------------------------------------------------
$iterRng = $d.Content.Duplicate
$iterRng.Find.Execute("iteration 161", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$iterEnd = $iterRng.End

# Force a run break between "iteration " and "161" by dropping a temporary
# bookmark at the split point, then overwrite the numeral.
$splitPoint = $d.Range($iterEnd - 3, $iterEnd - 3)
$d.Bookmarks.Add("_TempSplit", $splitPoint)

$numRng = $d.Range($iterEnd - 3, $iterEnd)
$numRng.Text = "140"

$d.Bookmarks("_TempSplit").Delete()

# Rewriting the run above makes the engine's re-run-ifier also coalesce the
# following, identically-formatted "." and " The result..." runs back into
# one; drop another transient bookmark at their former boundary so that
# split is preserved too.
$dotEnd = $iterEnd + 1
$splitPoint2 = $d.Range($dotEnd, $dotEnd)
$d.Bookmarks.Add("_TempSplit2", $splitPoint2)
$d.Bookmarks("_TempSplit2").Delete()

# "161" and "140" are the same length, so $iterEnd still marks the point
# right after the new "140" run - drop "_GoBack" there.
$goBackPoint = $d.Range($iterEnd, $iterEnd)
$d.Bookmarks.Add("_GoBack", $goBackPoint)
